$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 206.76923
$ws.Range("I33").Value = 151.6
$ws.Range("K33").Value = 151.6
$ws.Range("M33").Value = 77.40000000000001
$ws.Range("H132").Value = 2251.2144
$ws.Range("I132").Value = 1909.9565
$ws.Range("J132").Value = 3821
$ws.Range("K132").Value = 5729.8695
$ws.Range("L132").Value = 11463
$ws.Range("M132").Value = -3199.8695
$ws.Range("N132").Value = -16523
$ws.Range("H135").Value = 88236904
$ws.Range("I135").Value = 33335092
$ws.Range("J135").Value = 500000500
$ws.Range("K135").Value = 300015828
$ws.Range("L135").Value = 4500004500
$ws.Range("M135").Value = -300013293
$ws.Range("N135").Value = -4500009570
$ws.Range("H138").Value = 8822.371999999999
$ws.Range("J138").Value = 11101.579
$ws.Range("L138").Value = 33304.737
$ws.Range("N138").Value = -43584.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20620.611
$ws.Range("I32").Value = 20735.83
$ws.Range("J32").Value = 14514
$ws.Range("K32").Value = 20735.83
$ws.Range("L32").Value = 14514
$ws.Range("M32").Value = -20448.83
$ws.Range("N32").Value = -15088
$ws.Range("H110").Value = 1220.5385
$ws.Range("I110").Value = 1174.6
$ws.Range("J110").Value = 1373.6666
$ws.Range("K110").Value = 1174.6
$ws.Range("L110").Value = 1373.6666
$ws.Range("M110").Value = 870.4000000000001
$ws.Range("N110").Value = -5463.6666
$ws.Range("H122").Value = 1707.9131
$ws.Range("I122").Value = 1752
$ws.Range("J122").Value = 1498.5
$ws.Range("K122").Value = 5256
$ws.Range("L122").Value = 4495.5
$ws.Range("M122").Value = -2806
$ws.Range("N122").Value = -9395.5
$ws.Range("H132").Value = 12554.637
$ws.Range("I132").Value = 19846.5
$ws.Range("K132").Value = 59539.5
$ws.Range("M132").Value = -57009.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1099.8572
$ws.Range("I94").Value = 1216.5
$ws.Range("J94").Value = 400
$ws.Range("K94").Value = 1216.5
$ws.Range("L94").Value = 400
$ws.Range("M94").Value = -765.5
$ws.Range("N94").Value = -1302
$ws.Range("H134").Value = 38500.43
$ws.Range("I134").Value = 2571.3333
$ws.Range("K134").Value = 7713.999899999999
$ws.Range("M134").Value = -5178.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H50").Value = 28900
$ws.Range("J50").Value = 28900
$ws.Range("L50").Value = 28900
$ws.Range("N50").Value = -30150
$ws.Range("H51").Value = 39450
$ws.Range("J51").Value = 39450
$ws.Range("L51").Value = 39450
$ws.Range("N51").Value = -40922
$ws.Range("H60").Value = 16994.08
$ws.Range("I60").Value = 2000
$ws.Range("K60").Value = 2000
$ws.Range("M60").Value = -1489
$ws.Range("H61").Value = 39450
$ws.Range("J61").Value = 39450
$ws.Range("L61").Value = 39450
$ws.Range("N61").Value = -40146
$ws.Range("H122").Value = 14772
$ws.Range("I122").Value = 8678
$ws.Range("J122").Value = 30007
$ws.Range("K122").Value = 26034
$ws.Range("L122").Value = 90021
$ws.Range("M122").Value = -23584
$ws.Range("N122").Value = -94921
$ws.Range("H134").Value = 4215.421
$ws.Range("I134").Value = 3714.7
$ws.Range("J134").Value = 4771.778
$ws.Range("K134").Value = 11144.1
$ws.Range("L134").Value = 14315.334
$ws.Range("M134").Value = -8609.099999999999
$ws.Range("N134").Value = -19385.334
$ws.Range("H140").Value = 34850
$ws.Range("J140").Value = 34850
$ws.Range("L140").Value = 34850
$ws.Range("N140").Value = -45210

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2174.1428
$ws.Range("I70").Value = 1025
$ws.Range("J70").Value = 2444.5293
$ws.Range("K70").Value = 3075
$ws.Range("L70").Value = 7333.5879
$ws.Range("M70").Value = -2760
$ws.Range("N70").Value = -7963.5879
$ws.Range("H73").Value = 2174.1428
$ws.Range("I73").Value = 1025
$ws.Range("J73").Value = 2444.5293
$ws.Range("K73").Value = 3075
$ws.Range("L73").Value = 7333.5879
$ws.Range("M73").Value = -1983
$ws.Range("N73").Value = -9517.5879
$ws.Range("H129").Value = 1768.5
$ws.Range("I129").Value = 981.1111
$ws.Range("J129").Value = 2780.8572
$ws.Range("K129").Value = 2943.3333
$ws.Range("L129").Value = 8342.571599999999
$ws.Range("M129").Value = 2056.6667
$ws.Range("N129").Value = -18342.5716
$ws.Range("H131").Value = 994.6842
$ws.Range("I131").Value = 1086.5555
$ws.Range("J131").Value = 912
$ws.Range("K131").Value = 3259.6665
$ws.Range("L131").Value = 2736
$ws.Range("M131").Value = 1780.3335
$ws.Range("N131").Value = -12816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3173.3333
$ws.Range("J126").Value = 4200
$ws.Range("L126").Value = 12600
$ws.Range("N126").Value = -17540
$ws.Range("H132").Value = 30004.5
$ws.Range("I132").Value = 55012
$ws.Range("J132").Value = 21668.666
$ws.Range("K132").Value = 165036
$ws.Range("L132").Value = 65005.99800000001
$ws.Range("M132").Value = -162506
$ws.Range("N132").Value = -70065.99800000001
$ws.Range("H137").Value = 49766.668
$ws.Range("J137").Value = 49766.668
$ws.Range("L137").Value = 49766.668
$ws.Range("N137").Value = -59966.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2000
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -752
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H62").Value = 3937.75
$ws.Range("I62").Value = 4000.3333
$ws.Range("J62").Value = 3750
$ws.Range("K62").Value = 4000.3333
$ws.Range("L62").Value = 3750
$ws.Range("M62").Value = -3376.3333
$ws.Range("N62").Value = -4998
$ws.Range("H65").Value = 3937.75
$ws.Range("I65").Value = 4000.3333
$ws.Range("J65").Value = 3750
$ws.Range("K65").Value = 20001.6665
$ws.Range("L65").Value = 18750
$ws.Range("M65").Value = -16881.6665
$ws.Range("N65").Value = -24990
$ws.Range("H126").Value = 1809.5
$ws.Range("J126").Value = 1682
$ws.Range("L126").Value = 5046
$ws.Range("N126").Value = -9986
$ws.Range("H136").Value = 6008.9023
$ws.Range("I136").Value = 2710.125
$ws.Range("J136").Value = 10666
$ws.Range("K136").Value = 8130.375
$ws.Range("L136").Value = 31998
$ws.Range("M136").Value = -5580.375
$ws.Range("N136").Value = -37098
